$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 825
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
$ws.Range("H132").Value = 3311.7222
$ws.Range("I132").Value = 1660.5
$ws.Range("J132").Value = 5375.75
$ws.Range("K132").Value = 4981.5
$ws.Range("L132").Value = 16127.25
$ws.Range("M132").Value = -2451.5
$ws.Range("N132").Value = -21187.25
$ws.Range("H137").Value = 1247.5555
$ws.Range("I137").Value = 1103.8108
$ws.Range("J137").Value = 1333.3387
$ws.Range("K137").Value = 3311.4324
$ws.Range("L137").Value = 4000.0161
$ws.Range("M137").Value = -761.4323999999997
$ws.Range("N137").Value = -9100.016100000001
$ws.Range("H138").Value = 5559756
$ws.Range("I138").Value = 2580.4375
$ws.Range("J138").Value = 10005496
$ws.Range("K138").Value = 7741.3125
$ws.Range("L138").Value = 30016488
$ws.Range("M138").Value = -2601.3125
$ws.Range("N138").Value = -30026768
$ws.Range("H141").Value = 2471.1
$ws.Range("I141").Value = 2471.1
$ws.Range("K141").Value = 7413.299999999999
$ws.Range("M141").Value = -2233.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 187.6
$ws.Range("I4").Value = 187.6
$ws.Range("K4").Value = 187.6
$ws.Range("M4").Value = -71.59999999999999
$ws.Range("H32").Value = 3918.6191
$ws.Range("I32").Value = 2622.7322
$ws.Range("J32").Value = 14285.714
$ws.Range("K32").Value = 2622.7322
$ws.Range("L32").Value = 14285.714
$ws.Range("M32").Value = -2335.7322
$ws.Range("N32").Value = -14859.714
$ws.Range("H74").Value = 46686.91
$ws.Range("I74").Value = 84277.5
$ws.Range("K74").Value = 84277.5
$ws.Range("M74").Value = -83403.5
$ws.Range("H77").Value = 46686.91
$ws.Range("I77").Value = 84277.5
$ws.Range("K77").Value = 421387.5
$ws.Range("M77").Value = -417019.5
$ws.Range("H132").Value = 2299.5
$ws.Range("I132").Value = 1865.6316
$ws.Range("J132").Value = 3477.1428
$ws.Range("K132").Value = 5596.8948
$ws.Range("L132").Value = 10431.4284
$ws.Range("M132").Value = -3066.8948
$ws.Range("N132").Value = -15491.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5622.636
$ws.Range("I134").Value = 5729.56
$ws.Range("J134").Value = 5288.5
$ws.Range("K134").Value = 17188.68
$ws.Range("L134").Value = 15865.5
$ws.Range("M134").Value = -14653.68
$ws.Range("N134").Value = -20935.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 533.3570999999999
$ws.Range("I22").Value = 391.7
$ws.Range("J22").Value = 887.5
$ws.Range("K22").Value = 391.7
$ws.Range("L22").Value = 887.5
$ws.Range("M22").Value = -41.69999999999999
$ws.Range("N22").Value = -1587.5
$ws.Range("H31").Value = 8621870
$ws.Range("I31").Value = 717.41174
$ws.Range("J31").Value = 12196494
$ws.Range("K31").Value = 717.41174
$ws.Range("L31").Value = 12196494
$ws.Range("M31").Value = -422.41174
$ws.Range("N31").Value = -12197084
$ws.Range("H34").Value = 8621870
$ws.Range("I34").Value = 717.41174
$ws.Range("J34").Value = 12196494
$ws.Range("K34").Value = 717.41174
$ws.Range("L34").Value = 12196494
$ws.Range("M34").Value = -515.41174
$ws.Range("N34").Value = -12196898
$ws.Range("H58").Value = 12717
$ws.Range("I58").Value = 15853.286
$ws.Range("J58").Value = 1740
$ws.Range("K58").Value = 15853.286
$ws.Range("L58").Value = 1740
$ws.Range("M58").Value = -15650.286
$ws.Range("N58").Value = -2146
$ws.Range("H99").Value = 1851.9333
$ws.Range("I99").Value = 1842.0714
$ws.Range("J99").Value = 1990
$ws.Range("K99").Value = 1842.0714
$ws.Range("L99").Value = 1990
$ws.Range("M99").Value = -344.0714
$ws.Range("N99").Value = -4986
$ws.Range("H122").Value = 974
$ws.Range("I122").Value = 1008.8
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 3026.4
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -576.3999999999996
$ws.Range("N122").Value = -7300
$ws.Range("H126").Value = 1851.9333
$ws.Range("I126").Value = 1842.0714
$ws.Range("J126").Value = 1990
$ws.Range("K126").Value = 5526.2142
$ws.Range("L126").Value = 5970
$ws.Range("M126").Value = -3056.2142
$ws.Range("N126").Value = -10910
$ws.Range("H132").Value = 4722.857
$ws.Range("I132").Value = 4804.1
$ws.Range("J132").Value = 4519.75
$ws.Range("K132").Value = 14412.3
$ws.Range("L132").Value = 13559.25
$ws.Range("M132").Value = -11882.3
$ws.Range("N132").Value = -18619.25
$ws.Range("H134").Value = 3516.2307
$ws.Range("I134").Value = 3546.4546
$ws.Range("J134").Value = 3350
$ws.Range("K134").Value = 10639.3638
$ws.Range("L134").Value = 10050
$ws.Range("M134").Value = -8104.363799999999
$ws.Range("N134").Value = -15120
$ws.Range("H136").Value = 12717
$ws.Range("I136").Value = 15853.286
$ws.Range("J136").Value = 1740
$ws.Range("K136").Value = 47559.858
$ws.Range("L136").Value = 5220
$ws.Range("M136").Value = -45009.858
$ws.Range("N136").Value = -10320
$ws.Range("H140").Value = 69780
$ws.Range("J140").Value = 69780
$ws.Range("L140").Value = 69780
$ws.Range("N140").Value = -80140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 805.43335
$ws.Range("I122").Value = 516.75
$ws.Range("J122").Value = 1382.8
$ws.Range("K122").Value = 4650.75
$ws.Range("L122").Value = 12445.2
$ws.Range("M122").Value = -2200.75
$ws.Range("N122").Value = -17345.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 126114.125
$ws.Range("I122").Value = 200926.6
$ws.Range("J122").Value = 1426.6666
$ws.Range("K122").Value = 602779.8
$ws.Range("L122").Value = 4279.9998
$ws.Range("M122").Value = -600329.8
$ws.Range("N122").Value = -9179.9998
$ws.Range("H132").Value = 3648.0952
$ws.Range("I132").Value = 4973.143
$ws.Range("J132").Value = 2985.5715
$ws.Range("K132").Value = 14919.429
$ws.Range("L132").Value = 8956.7145
$ws.Range("M132").Value = -12389.429
$ws.Range("N132").Value = -14016.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2037.9354
$ws.Range("I136").Value = 1886.7826
$ws.Range("J136").Value = 2472.5
$ws.Range("K136").Value = 5660.3478
$ws.Range("L136").Value = 7417.5
$ws.Range("M136").Value = -3110.3478
$ws.Range("N136").Value = -12517.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 77111.44
$ws.Range("I122").Value = 2287.5557
$ws.Range("J122").Value = 173313.58
$ws.Range("K122").Value = 6862.6671
$ws.Range("L122").Value = 519940.74
$ws.Range("M122").Value = -4412.6671
$ws.Range("N122").Value = -524840.74
$ws.Range("H132").Value = 3291.5
$ws.Range("I132").Value = 3483.2
$ws.Range("J132").Value = 2932.0625
$ws.Range("K132").Value = 10449.6
$ws.Range("L132").Value = 8796.1875
$ws.Range("M132").Value = -7919.599999999999
$ws.Range("N132").Value = -13856.1875
